# Operate on the "Actual" worksheet (the daily attendance log for 2024-01-25).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Actual")

# Insert a new "Area" column before the existing "Fecha" column (C),
# shifting the old Fecha/Hora columns from C/D to D/E.
$ws.Range("C1").EntireColumn.Insert()

# Remove the last record (Mauricio Sanchez, row 5) entirely.
$ws.Range("A5").EntireRow.Delete()

# Header row.
$ws.Range("C1").Value = "Area"

# Row 2 -> replace the first record with a new one (Fabio Castada / Comercial).
# IDs and ISO dates must be entered as text (leading apostrophe) so Excel
# doesn't reinterpret them as numbers / date serials.
$ws.Range("A2").Value = "'2323232323"
$ws.Range("B2").Value = "Fabio Castada"
$ws.Range("C2").Value = "Comercial"
$ws.Range("D2").Value = "'2024-01-25"
$ws.Range("E2").Value = "22:4:13"

# Row 3 -> Julian Largo, Administrativa.
$ws.Range("A3").Value = "'1054398414"
$ws.Range("B3").Value = "Julian Largo"
$ws.Range("C3").Value = "Administrativa"
$ws.Range("D3").Value = "'2024-01-25"
$ws.Range("E3").Value = "22:5:54"

# Row 4 -> Julian Largo, Administrativa.
$ws.Range("A4").Value = "'1054398414"
$ws.Range("B4").Value = "Julian Largo"
$ws.Range("C4").Value = "Administrativa"
$ws.Range("D4").Value = "'2024-01-25"
$ws.Range("E4").Value = "22:7:54"
